$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel auto-coerces plain "123"/"45.6"-looking .Value assignments into
# numbers. All cells in this sheet must stay text (t="s", same style as
# before), so every write goes through a literal-text formula that is
# immediately collapsed back to a plain value via copy/paste-special
# (xlPasteValues = -4163). That keeps t="s" and the original style index
# intact instead of Excel inferring a numeric type / new number format.
function Set-TextCell($addr, $text) {
    $escaped = $text -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

# Row 1
Set-TextCell "A1" "{1 2018-11-15 16:26:10 +0800 CST 2018-11-15 16:26:10 +0800 CST <nil>}"
Set-TextCell "B1" "553.5"
Set-TextCell "C1" "9"
Set-TextCell "D1" "wabuguan"
Set-TextCell "E1" "online"
Set-TextCell "F1" "abc.txt"

# Row 2
Set-TextCell "A2" "{2 2018-11-15 16:26:51 +0800 CST 2018-11-15 16:34:44 +0800 CST <nil>}"
Set-TextCell "B2" "66.14"
Set-TextCell "C2" "1"
Set-TextCell "D2" "Alice"
Set-TextCell "E2" "offline"
Set-TextCell "F2" "image.jpg"

# Row 3
Set-TextCell "A3" "{3 2018-11-15 16:27:35 +0800 CST 2018-11-15 16:27:35 +0800 CST <nil>}"
Set-TextCell "B3" "43.24"
Set-TextCell "C3" "2"
Set-TextCell "D3" "Bob"
Set-TextCell "E3" "online"
Set-TextCell "F3" "bob.xlxs"

# Row 4
Set-TextCell "A4" "{4 2018-11-15 16:28:13 +0800 CST 2018-11-15 16:28:13 +0800 CST <nil>}"
Set-TextCell "B4" "99.66"
Set-TextCell "C4" "3"
Set-TextCell "D4" "Cat"
Set-TextCell "E4" "online"
Set-TextCell "F4" "cat.png"

# Row 5
Set-TextCell "A5" "{17 2018-11-15 17:57:14 +0800 CST 2018-11-15 17:57:14 +0800 CST <nil>}"
Set-TextCell "B5" "95.2"
Set-TextCell "C5" "5"
Set-TextCell "D5" "Elevant"
Set-TextCell "E5" "online"
Set-TextCell "F5" "Elevant.png"
